$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1800
$ws.Range("J29").Value = 2780
$ws.Range("L29").Value = 8340
$ws.Range("N29").Value = -8902
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968
$ws.Range("H58").Value = 723.86957
$ws.Range("J58").Value = 1024.2727
$ws.Range("L58").Value = 3072.8181
$ws.Range("N58").Value = -3372.8181
$ws.Range("I64").Value = 3677.6
$ws.Range("J64").Value = 3652
$ws.Range("K64").Value = 3677.6
$ws.Range("L64").Value = 3652
$ws.Range("M64").Value = -3429.6
$ws.Range("N64").Value = -4148
$ws.Range("I67").Value = 3677.6
$ws.Range("J67").Value = 3652
$ws.Range("K67").Value = 3677.6
$ws.Range("L67").Value = 3652
$ws.Range("M67").Value = -2819.6
$ws.Range("N67").Value = -5368
$ws.Range("H94").Value = 2833.3333
$ws.Range("I94").Value = 2500
$ws.Range("K94").Value = 2500
$ws.Range("M94").Value = -2049
$ws.Range("H98").Value = 3886
$ws.Range("I98").Value = 4857.846
$ws.Range("J98").Value = 1359.2
$ws.Range("K98").Value = 4857.846
$ws.Range("L98").Value = 1359.2
$ws.Range("M98").Value = -3359.846
$ws.Range("N98").Value = -4355.2
$ws.Range("H107").Value = 2045.92
$ws.Range("I107").Value = 1434.2632
$ws.Range("K107").Value = 1434.2632
$ws.Range("M107").Value = 485.7367999999999
$ws.Range("H116").Value = 2723.3684
$ws.Range("I116").Value = 2028.4445
$ws.Range("J116").Value = 3348.8
$ws.Range("K116").Value = 2028.4445
$ws.Range("L116").Value = 3348.8
$ws.Range("M116").Value = 1413.5555
$ws.Range("N116").Value = -10232.8
$ws.Range("H122").Value = 3886
$ws.Range("I122").Value = 4857.846
$ws.Range("J122").Value = 1359.2
$ws.Range("K122").Value = 14573.538
$ws.Range("L122").Value = 4077.6
$ws.Range("M122").Value = -12123.538
$ws.Range("N122").Value = -8977.6
$ws.Range("H138").Value = 1803.7
$ws.Range("I138").Value = 1381.2778
$ws.Range("J138").Value = 1896.4269
$ws.Range("K138").Value = 4143.8334
$ws.Range("L138").Value = 5689.280699999999
$ws.Range("M138").Value = 996.1665999999996
$ws.Range("N138").Value = -15969.2807
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 688.075
$ws.Range("I2").Value = 451.55554
$ws.Range("K2").Value = 451.55554
$ws.Range("M2").Value = -338.55554
$ws.Range("H32").Value = 7530.4316
$ws.Range("I32").Value = 6259.4165
$ws.Range("J32").Value = 13250
$ws.Range("K32").Value = 6259.4165
$ws.Range("L32").Value = 13250
$ws.Range("M32").Value = -5972.4165
$ws.Range("N32").Value = -13824
$ws.Range("H45").Value = 1172.5
$ws.Range("I45").Value = 1113.3334
$ws.Range("K45").Value = 1113.3334
$ws.Range("M45").Value = -736.3334
$ws.Range("H63").Value = 27029072
$ws.Range("I63").Value = 1942.2174
$ws.Range("J63").Value = 71430780
$ws.Range("K63").Value = 1942.2174
$ws.Range("L63").Value = 71430780
$ws.Range("M63").Value = -1256.2174
$ws.Range("N63").Value = -71432152
$ws.Range("H66").Value = 27029072
$ws.Range("I66").Value = 1942.2174
$ws.Range("J66").Value = 71430780
$ws.Range("K66").Value = 9711.087
$ws.Range("L66").Value = 357153900
$ws.Range("M66").Value = -6279.087
$ws.Range("N66").Value = -357160764
$ws.Range("H110").Value = 1443.2307
$ws.Range("I110").Value = 983.6
$ws.Range("K110").Value = 983.6
$ws.Range("M110").Value = 1061.4
$ws.Range("H116").Value = 688.075
$ws.Range("I116").Value = 451.55554
$ws.Range("K116").Value = 451.55554
$ws.Range("M116").Value = 1842.44446
$ws.Range("H132").Value = 2623.625
$ws.Range("I132").Value = 2152.3914
$ws.Range("K132").Value = 6457.174199999999
$ws.Range("M132").Value = -3927.174199999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 688.075
$ws.Range("I3").Value = 451.55554
$ws.Range("K3").Value = 451.55554
$ws.Range("M3").Value = -337.55554
$ws.Range("H54").Value = 6743.8335
$ws.Range("I54").Value = 1365.75
$ws.Range("J54").Value = 17500
$ws.Range("K54").Value = 1365.75
$ws.Range("L54").Value = 17500
$ws.Range("M54").Value = -881.75
$ws.Range("N54").Value = -18468
$ws.Range("H86").Value = 2652.9412
$ws.Range("I86").Value = 2710.15
$ws.Range("J86").Value = 2571.2144
$ws.Range("K86").Value = 2710.15
$ws.Range("L86").Value = 2571.2144
$ws.Range("M86").Value = -1587.15
$ws.Range("N86").Value = -4817.2144
$ws.Range("H89").Value = 2652.9412
$ws.Range("I89").Value = 2710.15
$ws.Range("J89").Value = 2571.2144
$ws.Range("K89").Value = 13550.75
$ws.Range("L89").Value = 12856.072
$ws.Range("M89").Value = -7934.75
$ws.Range("N89").Value = -24088.072
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 56040.4
$ws.Range("J92").Value = 56040.4
$ws.Range("L92").Value = 56040.4
$ws.Range("N92").Value = -61032.4
$ws.Range("H122").Value = 739.3333
$ws.Range("I122").Value = 714.61536
$ws.Range("K122").Value = 2143.84608
$ws.Range("M122").Value = 306.1539199999997
$ws.Range("H134").Value = 20001848
$ws.Range("I134").Value = 1734.4706
$ws.Range("J134").Value = 62502090
$ws.Range("K134").Value = 5203.4118
$ws.Range("L134").Value = 187506270
$ws.Range("M134").Value = -2668.4118
$ws.Range("N134").Value = -187511340
$ws.Range("H139").Value = 61759.668
$ws.Range("J139").Value = 61759.668
$ws.Range("L139").Value = 61759.668
$ws.Range("N139").Value = -72039.66800000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6603.4
$ws.Range("I56").Value = 6603.4
$ws.Range("K56").Value = 6603.4
$ws.Range("M56").Value = -6073.4
$ws.Range("H107").Value = 4379.2
$ws.Range("J107").Value = 5350
$ws.Range("L107").Value = 16050
$ws.Range("N107").Value = -19890
$ws.Range("H120").Value = 10232.75
$ws.Range("I120").Value = 4999
$ws.Range("J120").Value = 15466.5
$ws.Range("K120").Value = 14997
$ws.Range("L120").Value = 46399.5
$ws.Range("M120").Value = -10159
$ws.Range("N120").Value = -56075.5
$ws.Range("H131").Value = 20411006
$ws.Range("J131").Value = 3526.158
$ws.Range("L131").Value = 10578.474
$ws.Range("N131").Value = -20658.474
$ws.Range("H133").Value = 4126.8696
$ws.Range("J133").Value = 4383.222
$ws.Range("L133").Value = 13149.666
$ws.Range("N133").Value = -23269.666
$ws.Range("H137").Value = 22064462
$ws.Range("I137").Value = 68184070
$ws.Range("J137").Value = 7254.9565
$ws.Range("K137").Value = 204552210
$ws.Range("L137").Value = 21764.8695
$ws.Range("M137").Value = -204547110
$ws.Range("N137").Value = -31964.8695
$ws.Range("H141").Value = 66669700
$ws.Range("I141").Value = 90910740
$ws.Range("J141").Value = 6858.25
$ws.Range("K141").Value = 272732220
$ws.Range("L141").Value = 20574.75
$ws.Range("M141").Value = -272727040
$ws.Range("N141").Value = -30934.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7271.4287
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 5150
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 5150
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -7146
$ws.Range("H83").Value = 7271.4287
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 5150
$ws.Range("K83").Value = 100000
$ws.Range("L83").Value = 25750
$ws.Range("M83").Value = -95008
$ws.Range("N83").Value = -35734
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 1448.8
$ws.Range("J102").Value = 1612.6666
$ws.Range("L102").Value = 1612.6666
$ws.Range("N102").Value = -4856.6666
$ws.Range("H122").Value = 4766.9473
$ws.Range("J122").Value = 4163
$ws.Range("L122").Value = 12489
$ws.Range("N122").Value = -17389
$ws.Range("H126").Value = 1802.6086
$ws.Range("I126").Value = 1487.2
$ws.Range("K126").Value = 4461.6
$ws.Range("M126").Value = -1991.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5690.5386
$ws.Range("I40").Value = 2751.75
$ws.Range("K40").Value = 2751.75
$ws.Range("M40").Value = -2615.75
$ws.Range("H122").Value = 25001660
$ws.Range("I122").Value = 41668100
$ws.Range("J122").Value = 2001.25
$ws.Range("K122").Value = 125004300
$ws.Range("L122").Value = 6003.75
$ws.Range("M122").Value = -125001850
$ws.Range("N122").Value = -10903.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 541.1818
$ws.Range("J113").Value = 700.8
$ws.Range("L113").Value = 2102.4
$ws.Range("N113").Value = -6442.4
$ws.Range("H122").Value = 10001585
$ws.Range("I122").Value = 11365305
$ws.Range("J122").Value = 967.3333
$ws.Range("K122").Value = 34095915
$ws.Range("L122").Value = 2901.9999
$ws.Range("M122").Value = -34093465
$ws.Range("N122").Value = -7801.9999
